# "adjusted screws to match what was used"
# Update the screw/fastener quantities on the Mechanical sheet for a few
# part rows (M3x14, M3x16, M3 square nut, M2x8 self-tapping pan head) and
# leave the final view focused back on the Mechanical sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Mechanical")
$ws2 = $wb.Worksheets.Item("Electrical")

# --- Row 3 (M3x14): add Head Nod (F) and Head Tilt (G) screws, remove the
#     Chest Shoulder (J) ones that turned out not to be used ---
$ws1.Range("F3").Value = 2
$ws1.Range("G3").Value = 2
$ws1.Range("J3").ClearContents()

# --- Row 5 (M3x16): add Chest Shoulder (J) screws ---
$ws1.Range("J5").Value = 4

# --- Row 10 (M3 square nut): add Head Nod (F) and Head Tilt (G) nuts ---
$ws1.Range("F10").Value = 2
$ws1.Range("G10").Value = 2

# --- Row 11 (M2x8 self-tapping pan head): add Head Nod (F) and Head Tilt (G) screws ---
$ws1.Range("F11").Value = 4
$ws1.Range("G11").Value = 4

# --- Restore focus to the Mechanical tab with the cell that was last being
#     edited selected (mirrors the saved view state in the workbook) ---
$ws1.Activate()
$ws1.Range("I8").Select()
